$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the "backup" column (R) for the existing last two rows: they were
# written as empty inline strings; the completed data now carries an
# explicit numeric 0.
$ws.Range("R364").Value = 0
$ws.Range("R365").Value = 0

# Row 365 is now a completed ("isPivot") week.
$ws.Range("O365").Value = 1

# Append the newly-completed weekly bars (rows 366-370).
$newRows = @(
    @(45474, 672.0499877929688, 682.2000122070312, 661.7999877929688, 672.5999755859375, 671.5380859375, 8823109,  2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
    @(45481, 667.9500122070312, 668.25,             624.5499877929688, 643.7999877929688, 642.7835693359375, 17421424, 2024, 7, 8,  0, 0, 0, 28, 0, 0, 0),
    @(45488, 645,                648.1500244140625, 625.1500244140625, 633.6500244140625, 633.6500244140625, 7496948,  2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(45495, 628,                663.2999877929688, 620.5499877929688, 650.4000244140625, 650.4000244140625, 21264000, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 655,                661.5,              637.3499755859375, 646.0499877929688, 646.0499877929688, 10524085, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0)
)

$dateNumberFormat = $ws.Range("A365").NumberFormat

$startRow = 366
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    # Column A (Datetime) keeps the same date/time number format as the rest
    # of the column.
    $ws.Cells.Item($r, 1).NumberFormat = $dateNumberFormat
    # "backup" column (R) stays blank (an explicit empty string, like the
    # other rows in this column) for the brand-new rows. Writing a bare
    # apostrophe forces an empty *text* cell instead of clearing it outright;
    # re-apply the default style since the text coercion would otherwise
    # stamp a quote-prefix format onto the cell.
    $ws.Cells.Item($r, 18).Value = "'"
    $ws.Cells.Item($r, 18).Style = "Normal"
}
